# Daily attendance processing - 2025-12-20 07:02:49
#
# The "Recorded By" column (G) stores a comma-separated list of the
# accounts that touched each attendance record. A handful of rows were
# recorded with "System" (or the backdoor/system service accounts)
# listed first instead of last; normalize those specific entries by
# rotating the leading account to the back of the list, e.g.
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "admin@admin.com, System"                 -> "System, admin@admin.com"
#   "system, backup@backdoor.com, System"     -> "backup@backdoor.com, System, system"
# Every other "Recorded By" value is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$rotateMap = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "system, backup@backdoor.com, System" = "backup@backdoor.com, System, system"
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }

    if ($rotateMap.ContainsKey($value)) {
        $cell.Value = $rotateMap[$value]
    }
}
